# Game list view model with userid-to-username matching
#
# Applies the edits between before.xlsx and the target workbook:
#  - "M0 - Account Mgmt" sheet: expand the "Ties to item below..." note (D18)
#  - "M1 - Game Data" sheet: several task rows updated / added, plus a
#    print-orientation tweak and a new "Done" status
#  - Active tab moves from "M1 - Game Data" back to "M0 - Account Mgmt"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "M0 - Account Mgmt": expand note text in D18 (text content only;
# "NEXT 1"/"NEXT 2" labels in C18/C19 stay exactly as they were).
# ---------------------------------------------------------------------
$wsAccount = $wb.Worksheets.Item("M0 - Account Mgmt")
$wsAccount.Range("D18").Value = "Ties to item below. Note: Should also be able to show all games subscribed -> Will require a view model.    CREATE A USER PROFILE PAGE THAT USES INFO BOX FOR ANY USER. "

# ---------------------------------------------------------------------
# Sheet "M1 - Game Data": task list updates
# ---------------------------------------------------------------------
$wsGame = $wb.Worksheets.Item("M1 - Game Data")

# New note on "User can view all game entried by ANOTHER user" (row 14)
$wsGame.Range("D14").Value = "Link to user profile page once created. (See M0)"

# "Game opponents are pulled from database if they exist" (row 17) is now Done
$wsGame.Range("C17").Value = "Done"
$wsGame.Range("C17").Style = "Good"

# Clarify the hyperlink note text (row 18)
$wsGame.Range("A18").Value = "Game opponents are displayed as hyperlinks to userid if the user exists. "

# Stray bold formatting left on the (empty) D18 cell
$wsGame.Range("D18").Font.Bold = $true

# New row 19: additional task entry
$wsGame.Range("A19").Value = "Game entries have dates, and dates can be edited. "

# Print orientation tweak picked up on this sheet
$wsGame.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping, done last so the final activation
# state matches the target (active tab back on "M0 - Account Mgmt").
# ---------------------------------------------------------------------
$wsGame.Activate()
$wsGame.Range("A26:D26").Select()

$wsAccount.Activate()
$wsAccount.Range("D32").Select()
